$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.886.66'
$ws.Range('E2').Value = '  -1.07%  '

$ws.Range('D3').Value = '3.441.32'
$ws.Range('E3').Value = '  -0.09%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'584.15"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.62%  '

$ws.Range('D6').Value = "'173.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.87%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').Value = '3.438.96'
$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('E10').Value = '  -1.58%  '

$ws.Range('D11').Value = "'6.93"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.99%  '

$ws.Range('D12').Value = "'0.411"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.93%  '

$ws.Range('D13').Value = '4.037.21'
$ws.Range('E13').Value = '  -0.03%  '

$ws.Range('E14').Value = '  +1.74%  '

$ws.Range('D15').Value = "'28.86"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.78%  '

$ws.Range('D16').Value = '65.885.40'
$ws.Range('E16').Value = '  -1.01%  '

$ws.Range('D17').Value = "'0.0000171"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.58%  '

$ws.Range('D18').Value = '3.442.31'
$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('E19').Value = '  -1.09%  '

$ws.Range('D20').Value = "'13.75"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').Value = "'370.10"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.36%  '

$ws.Range('D22').Value = "'7.58"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '

$ws.Range('E23').Value = '  +2.09%  '

$ws.Range('D24').Value = "'0.998"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('E25').Value = '  +0.50%  '

$ws.Range('E26').Value = '  +3.27%  '

$ws.Range('E27').Value = '  -1.32%  '

$ws.Range('D28').Value = "'0.177"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.21%  '

$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('E30').Value = '  -1.52%  '

$ws.Range('D31').Value = "'23.57"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.25%  '

$ws.Range('E32').Value = '  -0.89%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('E34').Value = '  -4.69%  '

$ws.Range('D35').Value = "'6.98"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.59%  '

$ws.Range('E36').Value = '  +1.09%  '

$ws.Range('D37').Value = "'160.95"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.21%  '

$ws.Range('D38').Value = "'0.879"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.39%  '

$ws.Range('D39').Value = "'28.29"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.88%  '

$ws.Range('D40').Value = "'1.78"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.28%  '

$ws.Range('E41').Value = '  +0.22%  '

$ws.Range('D42').Value = '2.766.63'
$ws.Range('E42').Value = '  +2.83%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'4.45"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.53%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = "'6.45"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.00%  '

$ws.Range('E45').Value = '  -1.66%  '

$ws.Range('D46').Value = "'40.03"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.75%  '

$ws.Range('D47').Value = "'24.68"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.68%  '

$ws.Range('E48').Value = '  -1.04%  '

$ws.Range('D49').Value = "'323.07"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.99%  '

$ws.Range('E50').Value = '  +0.26%  '

$ws.Range('D51').Value = "'6.23"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.90%  '
